# Insert a new data row at row 316 (pushing the existing rows 316:446 down
# to 317:447) and populate it with the new weekly price-report entry.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Shift rows 316..446 down by one row, creating a blank row 316.
$ws.Rows.Item(316).Insert()

# Populate the newly inserted row 316 with the new "Ajo" (garlic) record.
$ws.Cells.Item(316, 1).Value = 4
$ws.Cells.Item(316, 2).Value = "Feria Lagunitas de Puerto Montt"
$ws.Cells.Item(316, 3).Value = "Los Lagos"
$ws.Cells.Item(316, 4).Value = 45027
$ws.Cells.Item(316, 5).Value = 10
$ws.Cells.Item(316, 6).Value = 100112003
$ws.Cells.Item(316, 7).Value = "Ajo"
$ws.Cells.Item(316, 8).Value = "Chino"
$ws.Cells.Item(316, 9).Value = "Primera"
$ws.Cells.Item(316, 10).Value = 240
$ws.Cells.Item(316, 11).Value = 18500
$ws.Cells.Item(316, 12).Value = 21000
$ws.Cells.Item(316, 13).Value = 19750
$ws.Cells.Item(316, 14).Value = "$/caja 10 kilos"
$ws.Cells.Item(316, 15).Value = "China"
$ws.Cells.Item(316, 16).Value = 1975
$ws.Cells.Item(316, 17).Value = 10
$ws.Cells.Item(316, 18).Value = "Hortaliza"
